# 1st major refactor pass complete.
#
# Two finished/obsolete tasks ("Engine refactoring (assume it is to be
# demoed" and "Finish off debug rendering (text at least)") sat at the top
# of the ToDo list in rows 2 and 3. They are removed here, which shifts
# every row below them up by two. The two existing cell-notes (on the
# "Compiled shaders" and "Error handling strategy..." rows) need to move
# up along with their rows so they stay attached to the same task text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- preserve the existing comments (they do not auto-follow a row delete) ---
# Capture full comment text before anything moves/deletes.
$commentB12Text = $ws.Range("B12").Comment.Text()
$commentB15Text = $ws.Range("B15").Comment.Text()

# Remove the old comments from their current (pre-delete) cells so they
# don't linger as orphans once the rows above them shift up.
$ws.Range("B12").Comment.Delete()
$ws.Range("B15").Comment.Delete()

# --- remove the two completed rows (old rows 2 and 3) ---
$ws.Range("A2:A3").EntireRow.Delete()

# After the delete, what used to be row 12 is now row 10, and what used to
# be row 15 is now row 13. Re-attach the comments there.
$ws.Range("B10").AddComment($commentB12Text)
$ws.Range("B13").AddComment($commentB15Text)

# --- restore the selection state (row 2 is no longer "done", row 3 is) ---
$ws.Range("A2:XFD3").Select()
$ws.Range("A3").Activate()
